$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 2.067959333333333
$ws.Range("N2").Value = 6.203878
$ws.Range("O2").Value = 0.03751906553627715
$ws.Range("P2").Value = 0.03751906553627715
$ws.Range("Q2").Value = 171.6151494756393
$ws.Range("R2").Value = 1544.536345280754
$ws.Range("S2").Value = 0.01684419991321273
$ws.Range("T2").Value = 0.01684419991321273
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("O3").Value = 0.1625861532004571
$ws.Range("P3").Value = 0.1625861532004571
$ws.Range("Q3").Value = 743.681821104712
$ws.Range("R3").Value = 6693.136389942409
$ws.Range("S3").Value = 0.07299312039050516
$ws.Range("T3").Value = 0.07299312039050515
$ws.Range("G4").Value = 82.98768099999999
$ws.Range("H4").Value = 248.963043
$ws.Range("I4").Value = 0.4489504115427952
$ws.Range("J4").Value = 0.4489504115427952
$ws.Range("O4").Value = 0.7998947812632657
$ws.Range("P4").Value = 0.7998947812632657
$ws.Range("Q4").Value = 3658.781488535448
$ws.Range("R4").Value = 32929.03339681904
$ws.Range("S4").Value = 0.3591130912390774
$ws.Range("T4").Value = 0.3591130912390774
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("M5").Value = 2.067959333333333
$ws.Range("N5").Value = 6.203878
$ws.Range("O5").Value = 0.03751906553627715
$ws.Range("P5").Value = 0.03751906553627715
$ws.Range("Q5").Value = 130.5721710240338
$ws.Range("R5").Value = 1175.149539216304
$ws.Range("S5").Value = 0.0128157902058829
$ws.Range("T5").Value = 0.0128157902058829
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("I6").Value = 0.3415807409566563
$ws.Range("J6").Value = 0.3415807409566563
$ws.Range("O6").Value = 0.1625861532004571
$ws.Range("P6").Value = 0.1625861532004571
$ws.Range("Q6").Value = 565.8250465034453
$ws.Range("R6").Value = 5092.425418531007
$ws.Range("S6").Value = 0.05553629867950457
$ws.Range("T6").Value = 0.05553629867950455
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("I7").Value = 0.3415807409566563
$ws.Range("J7").Value = 0.3415807409566563
$ws.Range("O7").Value = 0.7998947812632657
$ws.Range("P7").Value = 0.7998947812632657
$ws.Range("S7").Value = 0.2732286520712688
$ws.Range("T7").Value = 0.2732286520712688
$ws.Range("I8").Value = 0.2094688475005485
$ws.Range("J8").Value = 0.2094688475005485
$ws.Range("M8").Value = 2.067959333333333
$ws.Range("N8").Value = 6.203878
$ws.Range("O8").Value = 0.03751906553627715
$ws.Range("P8").Value = 0.03751906553627715
$ws.Range("Q8").Value = 80.07126544502533
$ws.Range("R8").Value = 720.641389005228
$ws.Range("S8").Value = 0.007859075417181522
$ws.Range("T8").Value = 0.007859075417181522
$ws.Range("I9").Value = 0.2094688475005485
$ws.Range("J9").Value = 0.2094688475005485
$ws.Range("O9").Value = 0.1625861532004571
$ws.Range("P9").Value = 0.1625861532004571
$ws.Range("R9").Value = 3122.847267134257
$ws.Range("S9").Value = 0.03405673413044735
$ws.Range("T9").Value = 0.03405673413044735
$ws.Range("I10").Value = 0.2094688475005485
$ws.Range("J10").Value = 0.2094688475005485
$ws.Range("O10").Value = 0.7998947812632657
$ws.Range("P10").Value = 0.7998947812632657
$ws.Range("S10").Value = 0.1675530379529196
$ws.Range("T10").Value = 0.1675530379529196
